# Updated cryptos list on Tue Feb  6 08:57:06 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "43.053.94"
$ws.Range("E2").Value = "  +0.03%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.333.21"
$ws.Range("E3").Value = "  +1.20%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.18%  "

# Row 5 - BNB
$ws.Range("D5").Value = "303.23"

# Row 6 - Solana
$ws.Range("D6").Value = "96.11"
$ws.Range("E6").Value = "  -1.22%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.503"
$ws.Range("E7").Value = "  -0.05%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.24%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -0.89%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "34.46"
$ws.Range("E10").Value = "  -2.91%  "

# Row 11 - Chainlink
$ws.Range("D11").Value = "19.15"
$ws.Range("E11").Value = "  +2.48%  "

# Row 12 - Dogecoin
$ws.Range("D12").Value = "0.0787"
$ws.Range("E12").Value = "  -0.28%  "

# Row 13 - TRON
$ws.Range("D13").Value = "0.123"
$ws.Range("E13").Value = "  +2.87%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "6.77"
$ws.Range("E14").Value = "  -1.95%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.695.28"
$ws.Range("E15").Value = "  +1.18%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.349.23"
$ws.Range("E16").Value = "  +0.95%  "

# Row 17 - Polygon
$ws.Range("D17").Value = "0.794"
$ws.Range("E17").Value = "  +1.38%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "43.007.52"
$ws.Range("E18").Value = "  +0.24%  "

# Row 19 - InternetComputer(DFINITY)
$ws.Range("D19").Value = "12.25"
$ws.Range("E19").Value = "  -2.77%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  +2.47%  "

# Row 21 - ShibaInu
$ws.Range("D21").Value = "0.0`u{2083}0894"
$ws.Range("E21").Value = "  -0.69%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "68.00"
$ws.Range("E22").Value = "  +0.59%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "236.73"
$ws.Range("E23").Value = "  -0.19%  "

# Row 24 - ImmutableX
$ws.Range("E24").Value = "  +3.60%  "

# Row 25 - was Dai, now PancakeSwap
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "2.43"
$ws.Range("E25").Value = "  +0.21%  "

# Row 26 - was PancakeSwap, now Dai
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.03%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "24.81"

# Row 28 - Toncoin
$ws.Range("E28").Value = "  -5.66%  "

# Row 29 - Cosmos
$ws.Range("D29").Value = "9.18"
$ws.Range("E29").Value = "  +1.19%  "

# Row 30 - InjectiveProtocol
$ws.Range("D30").Value = "32.06"
$ws.Range("E30").Value = "  -2.79%  "

# Row 31 - Monero
$ws.Range("D31").Value = "143.06"
$ws.Range("E31").Value = "  -14.02%  "

# Row 32 - FirstDigitalUSD
$ws.Range("E32").Value = "  +0.00%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "5.03"
$ws.Range("E33").Value = "  +0.72%  "

# Row 34 - Celestia
$ws.Range("D34").Value = "17.94"
$ws.Range("E34").Value = "  -1.28%  "

# Row 35 - was RenderToken, now Hedera
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.0704"
$ws.Range("E35").Value = "  +2.00%  "

# Row 36 - was Hedera, now RenderToken
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "4.43"
$ws.Range("E36").Value = "  -1.00%  "

# Row 37 - ARBITRUM
$ws.Range("E37").Value = "  +3.09%  "

# Row 38 - WEMIXToken
$ws.Range("D38").Value = "2.30"
$ws.Range("E38").Value = "  -2.14%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  +0.35%  "

# Row 40 - EnergySwap
$ws.Range("D40").Value = "22.39"
$ws.Range("E40").Value = "  +25.60%  "

# Row 41 - LidoDAOToken
$ws.Range("E41").Value = "  +0.07%  "

# Row 43 - Maker
$ws.Range("D43").Value = "1.935.70"
$ws.Range("E43").Value = "  -3.28%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  -0.06%  "

# Row 45 - FraxShare
$ws.Range("E45").Value = "  -2.79%  "

# Row 46 - ApeXProtocol
$ws.Range("E46").Value = "  -2.68%  "

# Row 47 - NEARProtocol
$ws.Range("E47").Value = "  -0.28%  "

# Row 48 - was RocketPoolETH, now HuobiToken
$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D48").Value = "2.88"
$ws.Range("E48").Value = "  +0.90%  "

# Row 49 - was HuobiToken, now RocketPoolETH
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.563.62"
$ws.Range("E49").Value = "  +1.24%  "

# Row 50 - MultiversX
$ws.Range("D50").Value = "53.79"
$ws.Range("E50").Value = "  +0.28%  "

# Row 51 - BitcoinSV
$ws.Range("D51").Value = "73.56"
$ws.Range("E51").Value = "  +2.35%  "
